# Apply the update described by the commit:
#  - bump the "Förändrad" (column C) date serial from 45184 to 45186 for every data row (2-51)
#  - add a friendly display-text argument (the "Beteckning", column A) to every HYPERLINK()
#    formula found in columns S, T, V, W, X, Y for the rows that contain them (2-8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = 2
$lastRow = $used.Row + $used.Rows.Count - 1

# Columns that may contain HYPERLINK formulas needing a friendly-name second argument.
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($row = $firstRow; $row -le $lastRow; $row++) {

    # --- Update column C (Förändrad) date value ---
    $cCell = $ws.Range("C$row")
    $cValue = $cCell.Value2
    if ($cValue -eq 45184) {
        $cCell.Value = 45186
    }

    # --- Update HYPERLINK formulas to include the display text ---
    $beteckning = $ws.Range("A$row").Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Range("$col$row")
        $formula = $cell.Formula
        if ($formula -and $formula.StartsWith("=HYPERLINK(")) {
            # Only patch formulas that don't already carry a second (friendly name) argument.
            if ($formula -notmatch '",\s*"') {
                $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
